# Append one new row (row 4) to Sheet1: the trip/shuttle record for
# "أحمد شريم" logged at 2025-05-01 17:14:32 (Arabic-formatted timestamp).
#
# Every existing cell in the sheet is stored as literal TEXT, including
# values that look numeric (e.g. C2="33", C3="222") and an intentionally
# empty string in column A. To reproduce that faithfully on write-back:
#   - the blank A-cell is entered with a leading apostrophe so Excel
#     commits an explicit empty TEXT value instead of clearing/omitting
#     the cell entirely (a plain "" assignment just deletes the cell), and
#   - the numeric-looking quantity ("22") is entered with a leading
#     apostrophe so Excel stores it as literal text rather than silently
#     converting it to a number.
# The apostrophe marks the cell with a "quote prefix" style; resetting the
# style back to "Normal" afterwards drops that marker while keeping the
# committed value as text, so the final cells come out identical in type
# and formatting to their neighbours in rows 1-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

$ws.Cells.Item($row, 1).Value = "'"                          # A4 -> "" (text)
$ws.Cells.Item($row, 2).Value = "أحمد شريم"                    # B4
$ws.Cells.Item($row, 3).Value = "'22"                         # C4 -> "22" (text)
$ws.Cells.Item($row, 4).Value = "الصمود"                      # D4
$ws.Cells.Item($row, 5).Value = "الرحلة 3"                     # E4
$ws.Cells.Item($row, 6).Value = "C1"                          # F4
$ws.Cells.Item($row, 7).Value = "UNICEF"                      # G4
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:١٤:٣٢ م"        # H4

# Drop the "quote prefix" formatting flag picked up from the leading
# apostrophes above so the new cells match the plain/default style used
# throughout the rest of the sheet.
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 3).Style = "Normal"
